# Fruta / hortaliza, semanal
#
# Two new weekly price rows were added to the daily logic subset for
# "Vega Central Mapocho de Santiago - Durazno". They belong right after the
# existing 2020-12-30 (serial 44166) batch, ahead of the 2021-02-13 (serial
# 44211) batch, which pushes every following row down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 392 (everything from the old row 392
# onward shifts down by two rows, to 394 onward).
$ws.Rows.Item(392).EntireRow.Insert()
$ws.Rows.Item(392).EntireRow.Insert()

# New row 392: Early Treat / Primera
$ws.Cells.Item(392, 1).Value = 9
$ws.Cells.Item(392, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(392, 3).Value = "Metropolitana"
$ws.Cells.Item(392, 4).Value = 44511
$ws.Cells.Item(392, 5).Value = 13
$ws.Cells.Item(392, 6).Value = "Fruta"
$ws.Cells.Item(392, 7).Value = 100103
$ws.Cells.Item(392, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(392, 9).Value = 100103004
$ws.Cells.Item(392, 10).Value = "Durazno"
$ws.Cells.Item(392, 11).Value = "Early Treat"
$ws.Cells.Item(392, 12).Value = "Primera"
$ws.Cells.Item(392, 13).Value = 400
$ws.Cells.Item(392, 14).Value = 11000
$ws.Cells.Item(392, 15).Value = 11000
$ws.Cells.Item(392, 16).Value = 11000
$ws.Cells.Item(392, 17).Value = "`$/bandeja 8 kilos empedrada"
$ws.Cells.Item(392, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(392, 19).Value = 1375
$ws.Cells.Item(392, 20).Value = 8

# New row 393: Florida King / Primera
$ws.Cells.Item(393, 1).Value = 9
$ws.Cells.Item(393, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(393, 3).Value = "Metropolitana"
$ws.Cells.Item(393, 4).Value = 44511
$ws.Cells.Item(393, 5).Value = 13
$ws.Cells.Item(393, 6).Value = "Fruta"
$ws.Cells.Item(393, 7).Value = 100103
$ws.Cells.Item(393, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(393, 9).Value = 100103004
$ws.Cells.Item(393, 10).Value = "Durazno"
$ws.Cells.Item(393, 11).Value = "Florida King"
$ws.Cells.Item(393, 12).Value = "Primera"
$ws.Cells.Item(393, 13).Value = 280
$ws.Cells.Item(393, 14).Value = 11000
$ws.Cells.Item(393, 15).Value = 11000
$ws.Cells.Item(393, 16).Value = 11000
$ws.Cells.Item(393, 17).Value = "`$/bandeja 8 kilos empedrada"
$ws.Cells.Item(393, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(393, 19).Value = 1375
$ws.Cells.Item(393, 20).Value = 8
